$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $val)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "64.738.44"
$ws.Range("E2").Value = "  -2.10%  "

Set-TextValue $ws "D3" "3.231.51"
$ws.Range("E3").Value = "  -1.17%  "

$ws.Range("E4").Value = "  +0.05%  "

Set-TextValue $ws "D5" "578.03"
$ws.Range("E5").Value = "  +0.44%  "

Set-TextValue $ws "D6" "175.26"
$ws.Range("E6").Value = "  -2.05%  "

Set-TextValue $ws "D7" "0.629"
$ws.Range("E7").Value = "  +0.80%  "

$ws.Range("E8").Value = "  +0.06%  "

Set-TextValue $ws "D9" "3.231.92"
$ws.Range("E9").Value = "  -1.08%  "

$ws.Range("E10").Value = "  -2.38%  "

$ws.Range("E11").Value = "  +0.71%  "

Set-TextValue $ws "D12" "0.392"
$ws.Range("E12").Value = "  -2.26%  "

Set-TextValue $ws "D13" "3.806.21"
$ws.Range("E13").Value = "  -0.65%  "

$ws.Range("E14").Value = "  -2.97%  "

Set-TextValue $ws "D15" "64.933.12"
$ws.Range("E15").Value = "  -1.87%  "

Set-TextValue $ws "D16" "25.69"
$ws.Range("E16").Value = "  -2.95%  "

Set-TextValue $ws "D17" "3.228.67"
$ws.Range("E17").Value = "  -1.18%  "

$ws.Range("E18").Value = "  -2.06%  "

Set-TextValue $ws "D19" "415.66"
$ws.Range("E19").Value = "  -3.98%  "

Set-TextValue $ws "D20" "5.37"
$ws.Range("E20").Value = "  -3.28%  "

Set-TextValue $ws "D21" "12.85"
$ws.Range("E21").Value = "  -2.26%  "

Set-TextValue $ws "D22" "7.20"
$ws.Range("E22").Value = "  -2.79%  "

Set-TextValue $ws "D23" "1.00"
$ws.Range("E23").Value = "  +0.02%  "

Set-TextValue $ws "D24" "70.26"
$ws.Range("E24").Value = "  -2.36%  "

Set-TextValue $ws "D25" "5.64"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("E26").Value = "  +3.68%  "

Set-TextValue $ws "D27" "0.498"
$ws.Range("E27").Value = "  -1.41%  "

$ws.Range("E28").Value = "  -2.33%  "

Set-TextValue $ws "D29" "9.21"
$ws.Range("E29").Value = "  +3.61%  "

$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("E31").Value = "  -2.84%  "

$ws.Range("E32").Value = "  -1.79%  "

$ws.Range("E33").Value = "  +0.00%  "

Set-TextValue $ws "D34" "5.04"
$ws.Range("E34").Value = "  -2.01%  "

$ws.Range("E35").Value = "  -2.29%  "

Set-TextValue $ws "D36" "1.15"
$ws.Range("E36").Value = "  -2.49%  "

Set-TextValue $ws "D37" "157.25"
$ws.Range("E37").Value = "  +0.30%  "

$ws.Range("E38").Value = "  -2.14%  "

Set-TextValue $ws "D39" "2.824.49"
$ws.Range("E39").Value = "  +1.85%  "

$ws.Range("E40").Value = "  -2.81%  "

Set-TextValue $ws "D41" "25.55"
$ws.Range("E41").Value = "  -3.50%  "

Set-TextValue $ws "D42" "4.22"
$ws.Range("E42").Value = "  -1.80%  "

Set-TextValue $ws "D43" "0.729"
$ws.Range("E43").Value = "  -5.71%  "

Set-TextValue $ws "D44" "39.21"
$ws.Range("E44").Value = "  -2.52%  "

Set-TextValue $ws "D45" "5.75"
$ws.Range("E45").Value = "  -4.33%  "

Set-TextValue $ws "D46" "0.0626"
$ws.Range("E46").Value = "  -4.43%  "

Set-TextValue $ws "D47" "305.01"
$ws.Range("E47").Value = "  -5.25%  "

# Rows 48 and 49 swap coin identity (B/C) and get new Price/Volume values.
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws "D48" "22.28"
$ws.Range("E48").Value = "  -4.46%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws "D49" "2.16"
$ws.Range("E49").Value = "  -6.36%  "

$ws.Range("E50").Value = "  -1.17%  "

$ws.Range("E51").Value = "  -0.46%  "
